$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.7287194209349384
$ws.Range("C2").Value = 9.226618575922256
$ws.Range("D2").Value = 0.1529057820181812
$ws.Range("E2").Value = 246.9852506941017
$ws.Range("G2").Value = 257.0934944729771

$ws.Range("B3").Value = 3.182878228561681
$ws.Range("C3").Value = 1.65323645889881
$ws.Range("D3").Value = 0.7127328510149897
$ws.Range("E3").Value = 0.4998867070740569
$ws.Range("G3").Value = 6.048734245549538

$ws.Range("B4").Value = 1.505614041169197
$ws.Range("C4").Value = 0.3375848360084654
$ws.Range("D4").Value = 157.8057217802531
$ws.Range("E4").Value = 0.4998867070740569
$ws.Range("G4").Value = 160.1488073645048
